$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the new "Yellow Area" grid-detection row (row 7) that was
# previously blank. This mirrors adding a new colour sample alongside
# Background / Wall / OutOfBounds.
$ws.Range("A7").Value = "Yellow Area"
$ws.Range("B7").Value = 60
$ws.Range("C7").Value = 79
$ws.Range("D7").Value = 83

# Recalculate so the dependent H/S/V formulas in E7:G7 pick up the new
# B7:D7 inputs.
$excel.Calculate() | Out-Null

# Leave the active selection where the author ended up after entering
# the new row of data.
$ws.Range("B9").Select() | Out-Null
